$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N2").Value = 3.45
$ws.Range("AG2").Value = 17.5
$ws.Range("AN2").Value = 55

$ws.Range("N4").Value = 1.36
$ws.Range("P4").Value = 1.36

$ws.Range("F8").Value = 1.98
$ws.Range("G8").Value = 2.54
$ws.Range("J8").Value = 3.4
$ws.Range("K8").Value = 7.4

$ws.Range("P10").Value = 3.3
$ws.Range("Q10").Value = 1.26

$ws.Range("F11").Value = 2.14
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 1.56
$ws.Range("I11").Value = 1.88
$ws.Range("J11").Value = 3.8
$ws.Range("P11").Value = 2.16

$ws.Range("P12").Value = 1.24
$ws.Range("Q12").Value = 1.02
